$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row => (C new value, E new value)
$updates = @{
    2   = @{ C = 766330;   E = 1429232685 }
    13  = @{ C = 187858;   E = 1168030789 }
    16  = @{ C = 10172;    E = 28006957 }
    27  = @{ C = 90066;    E = 442910133 }
    69  = @{ C = 17892;    E = 103956993 }
    78  = @{ C = 178445;   E = 892645000 }
    91  = @{ C = 18876;    E = 75305211 }
    121 = @{ C = 1306355;  E = 2275357834 }
    129 = @{ C = 633721;   E = 3433665537 }
    132 = @{ C = 585963;   E = 3471020874 }
    154 = @{ C = 18465;    E = 74174082 }
    237 = @{ C = 283323;   E = 1438433147 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
